$d = $word.ActiveDocument

# 1. "Phase II grant (up to $1,000,000)" -> "Phase II (up to $1,000,000)"
$d.Content.Find.Execute("Phase II grant (", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Phase II (", 2)

# 2. "Small businesses with Phase II grants are eligible" -> "Small businesses with Phase II funding are eligible"
$d.Content.Find.Execute("Small businesses with Phase II grants are eligible", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Small businesses with Phase II funding are eligible", 2)

# 3. "will know within three weeks if they meet the program's objectives" -> "will know within one month if they meet the program's objectives"
$d.Content.Find.Execute("will know within three weeks if they meet the program", $true, $false, $false, $false, $false,
                         $true, 1, $false, "will know within one month if they meet the program", 2)

# 4. "$1.75 million" -> "$2 million"
$d.Content.Find.Execute("up to $1.75 million", $true, $false, $false, $false, $false,
                         $true, 1, $false, "up to $2 million", 2)

# 5. "$8.1 billion" -> "$8.5 billion"
$d.Content.Find.Execute("a budget of about $8.1 billion", $true, $false, $false, $false, $false,
                         $true, 1, $false, "a budget of about $8.5 billion", 2)
